$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "simple-complex"

# Clear stale data that is not part of the new layout
$ws.Range("A2:A12").ClearContents()
$ws.Range("E1:E12").ClearContents()

# Write header row and data rows for the new layout
$ws.Cells.Item(1,1).Value = "question"
$ws.Cells.Item(1,2).Value = "input"
$ws.Cells.Item(1,3).Value = "analysis"
$ws.Cells.Item(1,4).Value = "response"
$ws.Cells.Item(1,6).Value = "analysis_3"
$ws.Cells.Item(1,7).Value = "response_3"
$ws.Cells.Item(1,9).Value = "analysis_2"
$ws.Cells.Item(1,10).Value = "response_2"
$ws.Cells.Item(2,2).Value = "Human wants to know the average pressure"
$ws.Cells.Item(2,3).Value = "I see this is too general question, and I do not have this information in my context so this type is complex"
$ws.Cells.Item(2,4).Value = "complex"
$ws.Cells.Item(2,6).Value = "I see this is too general question, Human has to specify a measurement system and a time range and I can not answer with the context"
$ws.Cells.Item(2,7).Value = "complex-incomplete"
$ws.Cells.Item(2,9).Value = "This is too general question, Human has to specify a measurement system and a time range"
$ws.Cells.Item(2,10).Value = "incomplete"
$ws.Cells.Item(3,2).Value = "Human wants to know information from a measurement system with tag 123123"
$ws.Cells.Item(3,3).Value = "This input is out of my knowledge, this type is complex"
$ws.Cells.Item(3,4).Value = "complex"
$ws.Cells.Item(3,6).Value = "This is too general because the measurement system has many information parameters and human has to specify wich parameters he wants and I can not answer with the context"
$ws.Cells.Item(3,7).Value = "complex-incomplete"
$ws.Cells.Item(3,9).Value = "This is too general because the measurement system has many information parameters and human has to specify wich parameters he wants"
$ws.Cells.Item(3,10).Value = "incomplete"
$ws.Cells.Item(4,2).Value = "Lauther is requesting information about a measurement system with the tag F980-40."
$ws.Cells.Item(4,3).Value = "This input is out of my knowledge, this type is complex"
$ws.Cells.Item(4,4).Value = "complex"
$ws.Cells.Item(4,6).Value = "This is too general, Lauther might specify wich information from measurement system wants also I can not answer with the context"
$ws.Cells.Item(4,7).Value = "complex-incomplete"
$ws.Cells.Item(4,9).Value = "This is too general, Lauther might specify wich information from measurement system wants"
$ws.Cells.Item(4,10).Value = "incomplete"
$ws.Cells.Item(5,2).Value = "The human wants to know the average pressure for the months of January and February for the measurement system with tag F980-40."
$ws.Cells.Item(5,3).Value = "It seems to be a detailed request but is out of my knowledge"
$ws.Cells.Item(5,4).Value = "complex"
$ws.Cells.Item(5,6).Value = "It seems to be a detailed request because user is giving a tag for measurement system, information parameter like pressure and a time range and I can not answer with the context"
$ws.Cells.Item(5,7).Value = "complex-complete"
$ws.Cells.Item(5,9).Value = "It seems to be a detailed request because user is giving a tag for measurement system, information parameter like pressure and a time range"
$ws.Cells.Item(5,10).Value = "complete"
$ws.Cells.Item(6,2).Value = "Andrew asks for a list of measurement systems in Spanish."
$ws.Cells.Item(6,3).Value = "This input is out of my knowledge, this type is complex"
$ws.Cells.Item(6,4).Value = "complex"
$ws.Cells.Item(6,6).Value = "It seems to be too general but user is asking for a list so in that order could be complete."
$ws.Cells.Item(6,7).Value = "complex-complete"
$ws.Cells.Item(6,9).Value = "It seems to be too general but user is asking for a list so in that order could be complete."
$ws.Cells.Item(6,10).Value = "complete"
$ws.Cells.Item(7,2).Value = "Alex wants is requesting the average gross flow for the months of October and December for the measurement system with tag AA9820-40."
$ws.Cells.Item(7,3).Value = "This input is out of my knowledge, this type is complex"
$ws.Cells.Item(7,4).Value = "complex"
$ws.Cells.Item(7,6).Value = "It seems to be a detailed request because Alex is giving a tag for measurement system, information parameter like pressure and a time range"
$ws.Cells.Item(7,7).Value = "complex-complete"
$ws.Cells.Item(7,9).Value = "It seems to be a detailed request because Andrew is giving a tag for measurement system, information parameter like pressure and a time range"
$ws.Cells.Item(7,10).Value = "complete"
$ws.Cells.Item(8,2).Value = "User is greeting AI and introduces himself as Lauther"
$ws.Cells.Item(8,3).Value = "It is a greeting-like, I can handle it"
$ws.Cells.Item(8,4).Value = "simple"
$ws.Cells.Item(8,6).Value = "It is a greeting-like"
$ws.Cells.Item(8,7).Value = "simple"
$ws.Cells.Item(8,9).Value = "It is a greeting-like"
$ws.Cells.Item(8,10).Value = "greeting-like"
$ws.Cells.Item(9,2).Value = "User is saying hello"
$ws.Cells.Item(9,3).Value = "User is greeting AI, I can say hello too"
$ws.Cells.Item(9,4).Value = "simple"
$ws.Cells.Item(9,6).Value = "User is greeting AI, there is no request"
$ws.Cells.Item(9,7).Value = "simple"
$ws.Cells.Item(9,9).Value = "User is greeting AI, there is no request"
$ws.Cells.Item(9,10).Value = "greeting-like"
$ws.Cells.Item(10,2).Value = "Alex wants is requesting the average viscosity for the months of October and December for the measurement system with tag AA9820-40."
$ws.Cells.Item(10,3).Value = "It seems to be a detailed request because Alex is giving a tag for measurement system, information parameter like pressure and a time range and I can not answer with the context"
$ws.Cells.Item(10,4).Value = "complex"
$ws.Cells.Item(10,6).Value = "It seems to be a detailed request because Alex is giving a tag for measurement system, information parameter like pressure and a time range and I can not answer with the context"
$ws.Cells.Item(10,7).Value = "complex-complete"
$ws.Cells.Item(10,9).Value = "It seems to be a detailed request because Andrew is giving a tag for measurement system, information parameter like pressure and a time range"
$ws.Cells.Item(10,10).Value = "complete"
$ws.Cells.Item(11,2).Value = "Andrew asks for a list of gas measurement systems"
$ws.Cells.Item(11,3).Value = "This input is out of my knowledge, this type is complex"
$ws.Cells.Item(11,4).Value = "complex"
$ws.Cells.Item(11,6).Value = "This is a complete request, Andrew is asking for a list of gas measurement systems and I can not answer with the context"
$ws.Cells.Item(11,7).Value = "complex-complete"
$ws.Cells.Item(11,9).Value = "This is a complete request, Andrew is asking for a list of gas measurement systems"
$ws.Cells.Item(11,10).Value = "complete"
$ws.Cells.Item(12,2).Value = "Andrew asks for a list of petrol measurement systems"
$ws.Cells.Item(12,3).Value = "This input is out of my knowledge, this type is complex"
$ws.Cells.Item(12,4).Value = "complex"
$ws.Cells.Item(12,6).Value = "This is a complete request, Andrew is asking for a list of petrol measurement systems and  I can not answer with the context"
$ws.Cells.Item(12,7).Value = "complex-complete"
$ws.Cells.Item(12,9).Value = "This is a complete request, Andrew is asking for a list of petrol measurement systems"
$ws.Cells.Item(12,10).Value = "complete"
$ws.Cells.Item(13,2).Value = "User is asking what AI is capable on"
$ws.Cells.Item(13,3).Value = "Based on my context It seems that is a question I can answer"
$ws.Cells.Item(13,4).Value = "simple"

# Update selection / view state
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.Zoom = 100
